# SummonTable.xlsx fix:
#  - ItemId column values switch from numeric item ids to string item codes
#  - Header "ItemId : Int" -> "ItemId : String"
#  - Right-align the ItemId data column
#  - Selection moves to E6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the ItemId column header (also renames the Table column since
#    B2 is the table's header cell).
$ws.Range("B2").Value = "ItemId : String"

# 2) Replace the literal ItemId values (numeric -> string code). All the
#    other B-column cells are formulas that copy these values down, so
#    they will pick up the new strings automatically on recalculation.
$idMap = @{
    "B3"  = "W0001"
    "B4"  = "W0002"
    "B5"  = "W0003"
    "B6"  = "W0005"
    "B7"  = "W0006"
    "B8"  = "W0007"
    "B9"  = "A0001"
    "B10" = "A0002"
    "B11" = "A0003"
    "B12" = "A0005"
    "B13" = "A0006"
    "B14" = "A0007"
    "B21" = "A0001"
    "B22" = "A0002"
    "B23" = "A0003"
    "B24" = "A0005"
    "B25" = "A0006"
    "B26" = "A0007"
}

foreach ($addr in $idMap.Keys) {
    $ws.Range($addr).Value = $idMap[$addr]
}

# 3) Right-align the whole ItemId data column (B3:B62) to match the new
#    string-based identifiers.
$ws.Range("B3:B62").HorizontalAlignment = -4152

# 4) Match the saved selection state.
$ws.Range("E6").Select()
